$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header row) - add P1, Q1 with header style (copy style from O1)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)  # xlPasteFormats

# Rows 2-25 - add P and Q columns with value 0
for ($r = 2; $r -le 25; $r++) {
    $ws.Cells.Item($r, 16).Value = 0
    $ws.Cells.Item($r, 17).Value = 0
}
